$d = $word.ActiveDocument

# --- Change 1: replace placeholder phone number with actual number ---
$d.Content.Find.Execute(
    "- Điện thoại: …………….",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "- Điện thoại: 0356.563.692",
    2) | Out-Null

# --- Change 2: expand permanent address with hamlet name ---
$d.Content.Find.Execute(
    "Canh Hiển, Vân Canh, Bình Định",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Thôn Hiển Đông, Canh Hiển, Vân Canh, Bình Định",
    2) | Out-Null

# --- Change 3: add company policy clause before working-hours text ---
$d.Content.Find.Execute(
    " Tối đa 8 giờ mỗi ngày và tối đa 40 giờ mỗi tuần",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    " theo quy định của công ty, tối đa 8 giờ mỗi ngày và tối đa 40 giờ mỗi tuần",
    2) | Out-Null
